# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values (per header in G1). Update the computed K
# values for each data row (rows 2-12) with the freshly regenerated
# numbers, replacing the previous "Strike#"-derived values.
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 2
